# Update "想去人数" (want-to-go count) values in the F column
# on sheet "展览" (Worksheets(1)) and "全部类型" (Worksheets(4)).
# Sheets "演出" and "本地生活" are unaffected by this change.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7349   # was 7299
$ws1.Range("F3").Value = 65     # was 63
$ws1.Range("F5").Value = 197    # was 179
$ws1.Range("F6").Value = 1112   # was 1109
$ws1.Range("F7").Value = 191    # was 188
$ws1.Range("F8").Value = 11     # was 10
$ws1.Range("F9").Value = 99     # was 92
$ws1.Range("F10").Value = 27    # was 25

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7349   # was 7299
$ws4.Range("F3").Value = 65     # was 63
$ws4.Range("F5").Value = 197    # was 179
$ws4.Range("F6").Value = 1112   # was 1109
$ws4.Range("F7").Value = 191    # was 188
$ws4.Range("F9").Value = 11     # was 10
$ws4.Range("F10").Value = 99    # was 92
$ws4.Range("F11").Value = 27    # was 25
